# "Add files via upload" - refresh the climate-data workbook:
#   - resize the saved window slightly
#   - clean up the "PRCP " header label (drop trailing space)
#   - append the new "US" national-average rows (2014-2022) under WYOMING
#   - leave the active selection on D3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- bookViews window size -------------------------------------------------
$win = $wb.Windows.Item(1)
$win.Width = 29400
$win.Height = 11700

# --- fix the PRCP column header (remove trailing space) --------------------
$ws.Range("D1").Value = "PRCP"

# --- append the new "US" rows (452-460) -------------------------------------
$usData = @(
    @(2014, 52.5283333333333,  2.5675),
    @(2015, 54.3825,           2.88166666666667),
    @(2016, 54.9016666666667,  2.6175),
    @(2017, 54.5383333333333,  2.69083333333333),
    @(2018, 53.5083333333333,  2.885),
    @(2019, 52.6641666666667,  2.90166666666667),
    @(2020, 54.3575,           2.53166666666667),
    @(2021, 53.3791666666667,  2.37083333333333),
    @(2022, 54.2125,           2.49666666666667)
)

$startRow = 452
for ($i = 0; $i -lt $usData.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $usData[$i][0]
    $ws.Cells.Item($r, 2).Value = "US"
    $ws.Cells.Item($r, 3).Value = $usData[$i][1]
    $ws.Cells.Item($r, 4).Value = $usData[$i][2]
}

# --- restore the selected cell ---------------------------------------------
$ws.Range("D3").Select()
